$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K (strikeout) column values, replacing the old "Strike#" data.
# Column G holds "K" per the header row; map of row -> new K value.
$newK = @{ 2=0; 3=1; 4=0; 5=1; 6=1; 7=0; 8=1; 9=1; 10=2; 11=1; 12=1; 13=2; 14=2; 15=0; 16=1; 17=1; 18=3; 19=0; 20=0; 21=1; 22=2; 23=1; 24=1; 25=0; 26=1; 27=1; 28=0; 29=2; 30=2; 31=2; 32=0; 33=3; 34=0; 35=1; 36=2; 37=1; 38=0; 39=1; 40=2; 41=0; 42=2; 43=1; 44=0; 45=0; 46=0; 47=1; 48=2; 49=1; 50=0; 51=0; 52=0; 53=0; 54=2; 55=0; 56=1; 57=1; 58=0; 59=2; 60=1; 61=2; 62=1; 63=0; 64=1; 65=1; 66=2; 67=1; 68=1; 69=2; 70=2; 71=0; 72=2; 73=0; 74=0; 75=0; 76=1; 77=1; 78=0; 80=2; 81=3; 82=3 }

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
